# Methodology slide (slide 6): reposition title/body placeholders and
# split the "GetInputFromExcelFile()" bullet into two runs
# ("GetInputFromExcelFile" + "( )").

$EMU_PER_POINT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# --- Title placeholder ("Methodology") ---------------------------------
$title = $s.Shapes.Item(1)

# Move the title box (size is unchanged).
$title.Left = 326571 / $EMU_PER_POINT
$title.Top  = 417927 / $EMU_PER_POINT

# Left-align the title paragraph.
$titleRange = $title.TextFrame.TextRange
$titleRange.ParagraphFormat.Alignment = 1   # ppAlignLeft

# --- Subtitle / body placeholder ----------------------------------------
$body = $s.Shapes.Item(2)

# Reposition + resize the body placeholder.
$body.Left   = 326571 / $EMU_PER_POINT
$body.Top    = 1301646 / $EMU_PER_POINT
$body.Width  = 11485984 / $EMU_PER_POINT
$body.Height = 5397734 / $EMU_PER_POINT

$bodyRange = $body.TextFrame.TextRange

# Find the "GetInputFromExcelFile()" bullet paragraph and split its single
# run into "GetInputFromExcelFile" + "( )".
for ($i = 1; $i -le $bodyRange.Paragraphs().Count; $i++) {
    $para = $bodyRange.Paragraphs($i, 1)
    if ($para.Text.TrimEnd() -eq "GetInputFromExcelFile()") {
        $parens = $para.Characters(22, 2)
        $parens.Text = "( )"
        break
    }
}
